$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update printer info
$ws.Range("C7").Value = "Nguyễn Nhật Lâm"
$ws.Range("C9").Value = "10000002"

# Update the overdue-loan data table (rows 13-17 hold the remaining 5 records,
# row 18 is cleared out since one record was removed from the list).
$ws.Range("B13").Value = "3"
$ws.Range("C13").Value = "2001170373"
$ws.Range("D13").Value = "Nguyễn Nhật Lâm"
$ws.Range("E13").Value = "26/01/2019"
$ws.Range("F13").Value = "3"
$ws.Range("G13").Value = "723"

$ws.Range("B14").Value = "4"
$ws.Range("C14").Value = "2001170018"
$ws.Range("D14").Value = "Nguyễn Nhật Lâm"
$ws.Range("E14").Value = "26/01/2019"
$ws.Range("F14").Value = "3"
$ws.Range("G14").Value = "723"

$ws.Range("B15").Value = "6"
$ws.Range("C15").Value = "2033207526"
$ws.Range("D15").Value = "Nguyễn Nhật Lâm"
$ws.Range("E15").Value = "26/01/2019"
$ws.Range("F15").Value = "1"
$ws.Range("G15").Value = "723"

$ws.Range("B16").Value = "7"
$ws.Range("C16").Value = "2033207526"
$ws.Range("D16").Value = "Nguyễn Nhật Lâm"
$ws.Range("E16").Value = "26/01/2019"
$ws.Range("F16").Value = "1"
$ws.Range("G16").Value = "723"

$ws.Range("B17").Value = "8"
$ws.Range("C17").Value = "2033207526"
$ws.Range("D17").Value = "Nguyễn Nhật Lâm"
$ws.Range("E17").Value = "26/01/2019"
$ws.Range("F17").Value = "1"
$ws.Range("G17").Value = "723"

# Row 18 no longer holds a record; clear all its cells.
$ws.Range("B18:G18").ClearContents()
